$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price/volume refresh). Values that are pure
# numeric-looking strings are written with a leading apostrophe so Excel
# keeps them as text (matching the source data, which stores these as
# plain strings, e.g. multi-dot prices like "29.208.18"), then the cell
# style is reset to Normal so no stray quote-prefix formatting is left
# behind on the cell.
$ws.Range("D2").Value = "29.208.18"
$ws.Range("D3").Value = "1.856.12"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'241.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.6995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.07793"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "'23.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("D11").Value = "'0.07803"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").Value = "1.860.54"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D14").Value = "'92.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "'0.6880"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "'6.548"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "'0.000008475"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").Value = "29.206.87"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'248.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").Value = "2.107.46"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'12.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'7.547"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'0.9997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'0.1507"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("D26").Value = "'161.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'8.862"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "'18.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("D29").Value = "'1.552"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("D30").Value = "'4.262"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").Value = "'4.210"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "'0.05232"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'0.7620"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").Value = "'1.170"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").Value = "'2.707"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'0.01861"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "1.227.07"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").Value = "'2.726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "'0.9005"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").Value = "'109.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "'0.9994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'5.538"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.66%  "
$ws.Range("D45").Value = "2.006.09"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").Value = "'65.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.74%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.553"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.5181"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "'1.750"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").Value = "'7.045"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
